# Generate Report for handoff
#
# - Status text "Handoff transform failed" -> "Ready for handoff" (Overview, zh-cn, de-de)
# - Add "Latest Handoff File" hyperlink cells (column C) on the zh-cn / de-de sheets
#   for the source-file row, pointing at the freshly produced xlf handoff files.
# - Populate the matching "Latest Handoff Datetime" (column D) for that row.
# - Flip "Handoff Reason" (column H) for that row from "Ignored" to "Include" now
#   that the file has actually been hand off.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# 1. Status message update - same shared text appears on all three sheets.
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsDeDe.Range("B2").Value = "Ready for handoff"

$commitSha = "9c8224f34e14126d086c77228f1b23ccd46c0890"
$repoBase  = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSha/e2e"

# 2. zh-cn sheet: new handoff file hyperlink + datetime + reason flip for row 2.
$zhFileName = "77cb3479-5038-4779-9011-473fd00eda82.b0fb71e3067d2bc5bb221eaacde3ee5f3962cc13.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "$repoBase/$zhFileName", $null, $null, $zhFileName)
$wsZhCn.Range("C2").Font.Underline = 2
$wsZhCn.Range("C2").Font.Color = 15570276

$wsZhCn.Range("D2").Value = "2016-01-11 03:50:59"
$wsZhCn.Range("H2").Value = "Include"

# 3. de-de sheet: new handoff file hyperlink + datetime + reason flip for row 2.
$deFileName = "77cb3479-5038-4779-9011-473fd00eda82.b0fb71e3067d2bc5bb221eaacde3ee5f3962cc13.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "$repoBase/$deFileName", $null, $null, $deFileName)
$wsDeDe.Range("C2").Font.Underline = 2
$wsDeDe.Range("C2").Font.Color = 15570276

$wsDeDe.Range("D2").Value = "2016-01-11 03:51:11"
$wsDeDe.Range("H2").Value = "Include"
